$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cost now precalculated in scenario: flip the yearly (2010-2050) cost
# multiplier cells for every ventilation-technology-efficiency-class row
# from 0 to 1.
$ws.Range("E2:AS8").Value = 1

# Match the author's resulting view state: select AM11 and scroll the
# viewport so column AJ is the left-most visible column.
$ws.Range("AM11").Select()
$excel.ActiveWindow.ScrollColumn = 36
$excel.ActiveWindow.ScrollRow = 1
